# Restore cell C10 ("Rules" sheet) to its prior value of 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1
